$p = $ppt.ActivePresentation

# --- Slide master & layouts: shorten the cached "datetimeFigureOut" field
#     text from the 4-digit year to a 2-digit year (8/15/2016 -> 8/15/16).
$m = $p.SlideMaster
$m.Shapes.Item(3).TextFrame.TextRange.Text = "8/15/16"

$cl = $m.CustomLayouts
$cl.Item(1).Shapes.Item(3).TextFrame.TextRange.Text = "8/15/16"
$cl.Item(2).Shapes.Item(3).TextFrame.TextRange.Text = "8/15/16"
$cl.Item(3).Shapes.Item(3).TextFrame.TextRange.Text = "8/15/16"
$cl.Item(4).Shapes.Item(4).TextFrame.TextRange.Text = "8/15/16"
$cl.Item(5).Shapes.Item(6).TextFrame.TextRange.Text = "8/15/16"
$cl.Item(6).Shapes.Item(2).TextFrame.TextRange.Text = "8/15/16"
$cl.Item(7).Shapes.Item(1).TextFrame.TextRange.Text = "8/15/16"
$cl.Item(8).Shapes.Item(4).TextFrame.TextRange.Text = "8/15/16"
$cl.Item(9).Shapes.Item(4).TextFrame.TextRange.Text = "8/15/16"
$cl.Item(10).Shapes.Item(3).TextFrame.TextRange.Text = "8/15/16"
$cl.Item(11).Shapes.Item(3).TextFrame.TextRange.Text = "8/15/16"

# --- Slide 1: renumber the four step callouts (1,2,3,4) -> (2,3,4,1)
#     and nudge the third bubble ("Oval 42") to the right.
$s = $p.Slides.Item(1)
$s.Shapes.Item(17).TextFrame.TextRange.Text = "2"
$s.Shapes.Item(18).TextFrame.TextRange.Text = "3"
$s.Shapes.Item(19).TextFrame.TextRange.Text = "4"
$s.Shapes.Item(19).Left = 798.3375590551182
$s.Shapes.Item(20).TextFrame.TextRange.Text = "1"
